$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# Text format first, so Excel stores them as strings (matching the source
# inlineStr cells) instead of auto-converting to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.167.24'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = '2.267.53'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '306.01'
$ws.Range("D6").Value = '97.25'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("D7").Value = '0.526'
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").Value = '35.15'
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").Value = '0.0789'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '6.88'
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").Value = '2.617.93'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").Value = '2.260.61'
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("D17").Value = '0.790'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").Value = '42.051.65'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").Value = '12.26'
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '236.94'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("D25").Value = '2.59'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '23.51'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = '37.50'
$ws.Range("E28").Value = '  +3.97%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '2.12'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = '162.34'
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").Value = '5.25'
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '3.15'
$ws.Range("E34").Value = '  +2.06%  '
$ws.Range("D35").Value = '17.61'
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("D36").Value = '0.0735'
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("D43").Value = '19.08'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("D44").Value = '1.948.18'
$ws.Range("E44").Value = '  -3.09%  '
$ws.Range("D45").Value = '0.0280'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").Value = '9.92'
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("E47").Value = '  -2.41%  '
$ws.Range("D48").Value = '53.85'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '2.490.55'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("D50").Value = '92.15'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = '71.60'
$ws.Range("E51").Value = '  -1.71%  '
